$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column C. This shifts the old C column (ACTION /
#    markup.setValue / MarkupValue / "10") to column D, carrying formatting and
#    the merged-cell ranges (B1:C1 -> B1:D1, etc.) along with it.
$ws.Columns("C:C").Insert()

# 2. Populate the new condition column (C) for the decision table:
#    row 7 "CONDITION" header, row 8 the "$journey:Journey" binding, row 9 the
#    new condition formula text, row 10 the column name.
$ws.Range("C7").Value = "CONDITION"
$ws.Range("C8").Value = "`$journey:Journey"
$ws.Range("C9").Value = "`$journey.getArrLocation().getCountryCode() in (`$param)"
$ws.Range("C10").Value = "JourneyLocationCountryCode"

# C11 (first rule row) stays empty for the new condition column, but needs the
# same formatting as its neighbour B11.
$ws.Range("B11").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# 3. Add a new rule row (row 12) "second_rule_on_markup": condition "FR" for the
#    new ArrLocation column, and MarkupValue action "20". Column B (DeptLocation
#    condition) is left blank for this rule.
$ws.Range("A12").Value = "second_rule_on_markup"
$ws.Range("C12").Value = """FR"""
$ws.Range("D12").Value = """20"""

# Copy formatting from row 11 into row 12 so styles line up. Note: row 12's
# last column (D) uses the plain "rule value" style (like B11/C11), not the
# highlighted style of D11, so copy A11:C11 -> A12:C12 and B11 -> D12.
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("D12").PasteSpecial(-4122)

# Re-apply values after the format paste (PasteSpecial(formats) shouldn't touch
# values, but make sure content is exactly as intended).
$ws.Range("A12").Value = "second_rule_on_markup"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = """FR"""
$ws.Range("D12").Value = """20"""

# 4. Resize columns to match the new layout (A/B narrower, C new, D keeps the
#    old column-C width). ColumnWidth is in "characters" and Excel snaps it to
#    a pixel grid, so these inputs are chosen to land as close as possible to
#    the intended 22.43 / 48.93 / 49.1 / 45.91 character widths.
$ws.Columns("A:A").ColumnWidth = 21.666666666666664
$ws.Columns("B:B").ColumnWidth = 48.16666666666667
$ws.Columns("C:C").ColumnWidth = 48.33333333333333
$ws.Columns("D:D").ColumnWidth = 45.0

# 5. Update the active selection to match the edited workbook (D11).
$ws.Range("D11").Select()
